$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 382; existing rows 382..466 shift down to 383..467
$ws.Rows.Item(382).Insert()

# Populate the newly inserted row 382 with the new record's data
$ws.Range("A382").Value = 6
$ws.Range("B382").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C382").Value = "Metropolitana"
$ws.Range("D382").Value = 45211
$ws.Range("E382").Value = 13
$ws.Range("F382").Value = 100112026
$ws.Range("G382").Value = "Haba"
$ws.Range("H382").Value = "Sin especificar"
$ws.Range("I382").Value = "Primera"
$ws.Range("J382").Value = 550
$ws.Range("K382").Value = 5000
$ws.Range("L382").Value = 6000
$ws.Range("M382").Value = 5455
$ws.Range("N382").Value = "`$/saco 25 kilos"
$ws.Range("O382").Value = "Región Metropolitana"
$ws.Range("P382").Value = 218
$ws.Range("Q382").Value = 25
$ws.Range("R382").Value = "Hortaliza"
